$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old H1 "derivatizable" header cell/column — clearing its
#    content drops the now-unused shared string and shrinks the used range
#    from A1:H73 down to A1:G73.
$ws.Range("H1").ClearContents()

# 2. Populate the new "detectable" flag column (G) for rows 2-12 with 1.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 7).Value = 1
}

# 3. Highlight the "Fatty alcohols" sub_class label in red.
$ws.Range("B33").Font.Color = 255

# 4. Update the active selection to G13 (also clears the stale
#    topLeftCell="B1" scroll position on the sheet view).
[void]$ws.Range("G13").Select()
